$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 holds a date serial number; it moves forward one month:
# 45406 (2024-04-24) -> 45436 (2024-05-24)
$ws.Range("A1").Value = 45436

# Update the two price cells in the price list.
$ws.Range("D29").Value = 520.458
$ws.Range("D30").Value = 353.073
